$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.364.37"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "3.525.50"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.10"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.66"
$ws.Range("E6").Value = "  +1.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  +1.85%  "

$ws.Range("E9").Value = "  +6.57%  "

$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").Value = "4.139.04"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.74"
$ws.Range("E14").Value = "  +1.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000182"
$ws.Range("E15").Value = "  +0.78%  "

$ws.Range("D16").Value = "67.279.04"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").Value = "3.521.12"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.26"
$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "397.45"
$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.99"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.52"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.541"
$ws.Range("E23").Value = "  +1.74%  "

$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("E25").Value = "  -4.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.34"
$ws.Range("E26").Value = "  +1.26%  "

$ws.Range("E27").Value = "  -0.66%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.30"
$ws.Range("E29").Value = "  -1.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("E30").Value = "  -1.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.09"
$ws.Range("E31").Value = "  +0.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.14"
$ws.Range("E32").Value = "  +2.07%  "

$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("E34").Value = "  +2.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.89"
$ws.Range("E35").Value = "  +1.44%  "

$ws.Range("E36").Value = "  -1.62%  "

$ws.Range("E37").Value = "  -1.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.93"
$ws.Range("E38").Value = "  +3.01%  "

$ws.Range("E39").Value = "  +1.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0748"
$ws.Range("E40").Value = "  -0.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.46"
$ws.Range("E41").Value = "  +1.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.52"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("E43").Value = "  +1.87%  "

$ws.Range("D44").Value = "2.804.68"
$ws.Range("E44").Value = "  -0.90%  "

$ws.Range("E45").Value = "  -1.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0310"
$ws.Range("E46").Value = "  -1.94%  "

$ws.Range("E47").Value = "  -3.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.10"
$ws.Range("E48").Value = "  +0.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.91"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("E50").Value = "  +0.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.852"
$ws.Range("E51").Value = "  -0.19%  "
